# Legislator property workbook update (#5: insurance, claim, debt, investment done)
# Adds proper header rows + legislator/source metadata columns to the
# 保險 (insurance), 債權 (claim), 債務 (debt) and 事業投資 (investment) sheets,
# matching the schema already used by the other sheets in the workbook, and
# fixes the mis-tagged property_category ("otherbonds" -> "antique") on the
# 具有相當價值之財產 (valuable property) sheet.

$wb = $excel.ActiveWorkbook

function Set-HeaderCell($ws, $addr, $text) {
    $r = $ws.Range($addr)
    $r.Value = $text
    $r.Font.Bold = $true
    $r.HorizontalAlignment = -4108
    $r.VerticalAlignment = -4160
    $r.Borders.LineStyle = 1
}

# ---------------------------------------------------------------------
# 具有相當價值之財產 (sheet index 5) - retag property_category for the
# three antique/jewellery rows from the stray "otherbonds" value to the
# correct "antique" value. (This also fixes up the shared-string text
# used elsewhere.)
# ---------------------------------------------------------------------
$wsValuable = $wb.Worksheets.Item(5)
$wsValuable.Range("F2").Value = "antique"
$wsValuable.Range("F3").Value = "antique"
$wsValuable.Range("F4").Value = "antique"

# ---------------------------------------------------------------------
# 保險 (insurance) - sheet index 6
# ---------------------------------------------------------------------
$wsIns = $wb.Worksheets.Item(6)

Set-HeaderCell $wsIns "B1" "company"
Set-HeaderCell $wsIns "C1" "name"
Set-HeaderCell $wsIns "D1" "owner"
Set-HeaderCell $wsIns "E1" "property_category"
Set-HeaderCell $wsIns "F1" "category"
Set-HeaderCell $wsIns "G1" "date"
Set-HeaderCell $wsIns "H1" "legislator_name"
Set-HeaderCell $wsIns "I1" "legislator_id"
Set-HeaderCell $wsIns "J1" "source_file"
Set-HeaderCell $wsIns "K1" "index"

$wsIns.Range("B2").Value = "富邦人壽"
$wsIns.Range("C2").Value = "投資型壽險"
$wsIns.Range("D2").Value = "姚文智"
$wsIns.Range("E2").Value = "insurance"
$wsIns.Range("F2").Value = "normal"
$wsIns.Range("G2").Value = "2013-12-31"
$wsIns.Range("H2").Value = "姚文智"
$wsIns.Range("I2").Value = 1745
$wsIns.Range("J2").Value = "tmpc2191"
$wsIns.Range("K2").Value = 86

$wsIns.Range("B3").Value = "南山人壽"
$wsIns.Range("C3").Value = "投資型壽險"
$wsIns.Range("D3").Value = "潘瓊琪"
$wsIns.Range("E3").Value = "insurance"
$wsIns.Range("F3").Value = "normal"
$wsIns.Range("G3").Value = "2013-12-31"
$wsIns.Range("H3").Value = "姚文智"
$wsIns.Range("I3").Value = 1745
$wsIns.Range("J3").Value = "tmpc2191"
$wsIns.Range("K3").Value = 87

# ---------------------------------------------------------------------
# 債權 (claim) - sheet index 7
# ---------------------------------------------------------------------
$wsClaim = $wb.Worksheets.Item(7)

Set-HeaderCell $wsClaim "B1" "species"
Set-HeaderCell $wsClaim "C1" "owner"
Set-HeaderCell $wsClaim "D1" "debtor"
Set-HeaderCell $wsClaim "E1" "total"
Set-HeaderCell $wsClaim "F1" "register_date"
Set-HeaderCell $wsClaim "G1" "register_reason"
Set-HeaderCell $wsClaim "H1" "property_category"
Set-HeaderCell $wsClaim "I1" "category"
Set-HeaderCell $wsClaim "J1" "date"
Set-HeaderCell $wsClaim "K1" "legislator_name"
Set-HeaderCell $wsClaim "L1" "legislator_id"
Set-HeaderCell $wsClaim "M1" "source_file"
Set-HeaderCell $wsClaim "N1" "index"

$wsClaim.Range("B2").Value = "借款"
$wsClaim.Range("C2").Value = "潘瓊琪"
$wsClaim.Range("D2").Value = "創意企業有限公司新北市新店區民權路88號4F"
$wsClaim.Range("E2").Value = 2000000
$wsClaim.Range("F2").Value = "102年08月01日"
$wsClaim.Range("G2").Value = "借款"
$wsClaim.Range("H2").Value = "claim"
$wsClaim.Range("I2").Value = "normal"
$wsClaim.Range("J2").Value = "2013-12-31"
$wsClaim.Range("K2").Value = "姚文智"
$wsClaim.Range("L2").Value = 1745
$wsClaim.Range("M2").Value = "tmpc2191"
$wsClaim.Range("N2").Value = 92

# ---------------------------------------------------------------------
# 債務 (debt) - sheet index 8
# ---------------------------------------------------------------------
$wsDebt = $wb.Worksheets.Item(8)

Set-HeaderCell $wsDebt "B1" "species"
Set-HeaderCell $wsDebt "C1" "debtor"
Set-HeaderCell $wsDebt "D1" "owner"
Set-HeaderCell $wsDebt "E1" "total"
Set-HeaderCell $wsDebt "F1" "register_date"
Set-HeaderCell $wsDebt "G1" "register_reason"
Set-HeaderCell $wsDebt "H1" "property_category"
Set-HeaderCell $wsDebt "I1" "category"
Set-HeaderCell $wsDebt "J1" "date"
Set-HeaderCell $wsDebt "K1" "legislator_name"
Set-HeaderCell $wsDebt "L1" "legislator_id"
Set-HeaderCell $wsDebt "M1" "source_file"
Set-HeaderCell $wsDebt "N1" "index"

$wsDebt.Range("B2").Value = "房屋貸款"
$wsDebt.Range("C2").Value = "潘瓊琪"
$wsDebt.Range("D2").Value = "花旗(台灣)商業銀行臺北市信義區松智路1號"
$wsDebt.Range("E2").Value = 8070000
$wsDebt.Range("F2").Value = "93年10月01日"
$wsDebt.Range("G2").Value = "轉貸自合庫"
$wsDebt.Range("H2").Value = "debt"
$wsDebt.Range("I2").Value = "normal"
$wsDebt.Range("J2").Value = "2013-12-31"
$wsDebt.Range("K2").Value = "姚文智"
$wsDebt.Range("L2").Value = 1745
$wsDebt.Range("M2").Value = "tmpc2191"
$wsDebt.Range("N2").Value = 97

# ---------------------------------------------------------------------
# 事業投資 (investment) - sheet index 9
# ---------------------------------------------------------------------
$wsInv = $wb.Worksheets.Item(9)

Set-HeaderCell $wsInv "B1" "owner"
Set-HeaderCell $wsInv "C1" "company"
Set-HeaderCell $wsInv "D1" "address"
Set-HeaderCell $wsInv "E1" "total"
Set-HeaderCell $wsInv "F1" "register_date"
Set-HeaderCell $wsInv "G1" "register_reason"
Set-HeaderCell $wsInv "H1" "property_category"
Set-HeaderCell $wsInv "I1" "category"
Set-HeaderCell $wsInv "J1" "date"
Set-HeaderCell $wsInv "K1" "legislator_name"
Set-HeaderCell $wsInv "L1" "legislator_id"
Set-HeaderCell $wsInv "M1" "source_file"
Set-HeaderCell $wsInv "N1" "index"

$wsInv.Range("B2").Value = "潘瓊琪"
$wsInv.Range("C2").Value = "創意企業有限公司"
$wsInv.Range("D2").Value = "新北市民權路88號4F"
$wsInv.Range("E2").Value = 2500000
$wsInv.Range("F2").Value = "93年10月01日"
$wsInv.Range("G2").Value = "投資"
$wsInv.Range("H2").Value = "investment"
$wsInv.Range("I2").Value = "normal"
$wsInv.Range("J2").Value = "2013-12-31"
$wsInv.Range("K2").Value = "姚文智"
$wsInv.Range("L2").Value = 1745
$wsInv.Range("M2").Value = "tmpc2191"
$wsInv.Range("N2").Value = 102
